$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename "Shrub" vegetation treatment labels to "CSS" in column A (labels column)
$ws.Range("A2").Value = "0 x CSS"
$ws.Range("A3").Value = "3 x CSS"
$ws.Range("A4").Value = "5 x CSS"
$ws.Range("A5").Value = "6 x CSS"
